$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2316026
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 2382192.2
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 7146576.600000001
$ws.Range("M17").Value = -432
$ws.Range("N17").Value = -7146912.600000001

$ws.Range("H28").Value = 3100.3635
$ws.Range("I28").Value = 3461
$ws.Range("J28").Value = 2799.8333
$ws.Range("K28").Value = 3461
$ws.Range("L28").Value = 2799.8333
$ws.Range("M28").Value = -2976
$ws.Range("N28").Value = -3769.8333

$ws.Range("H40").Value = 1999.1111
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 1998.6666
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 1998.6666
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2348.6666

$ws.Range("H51").Value = 3376
$ws.Range("I51").Value = 2590
$ws.Range("K51").Value = 2590
$ws.Range("M51").Value = -2106

$ws.Range("H87").Value = 26924
$ws.Range("J87").Value = 26924
$ws.Range("L87").Value = 26924
$ws.Range("N87").Value = -29420

$ws.Range("H90").Value = 26924
$ws.Range("J90").Value = 26924
$ws.Range("L90").Value = 80772
$ws.Range("N90").Value = -93252

$ws.Range("H101").Value = 1661.8889
$ws.Range("I101").Value = 967.4286
$ws.Range("J101").Value = 4092.5
$ws.Range("K101").Value = 2902.2858
$ws.Range("L101").Value = 12277.5
$ws.Range("M101").Value = -1280.2858
$ws.Range("N101").Value = -15521.5

$ws.Range("H129").Value = 1772.238
$ws.Range("I129").Value = 354.25
$ws.Range("K129").Value = 1062.75
$ws.Range("M129").Value = 3937.25

$ws.Range("H139").Value = 36211.668
$ws.Range("J139").Value = 43963
$ws.Range("L139").Value = 43963
$ws.Range("N139").Value = -54243

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12399.116
$ws.Range("I32").Value = 11953.2
$ws.Range("K32").Value = 11953.2
$ws.Range("M32").Value = -11666.2

$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 3000
$ws.Range("K63").Value = 3000
$ws.Range("M63").Value = -2314

$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 3000
$ws.Range("K66").Value = 15000
$ws.Range("M66").Value = -11568

$ws.Range("H88").Value = 2959.9
$ws.Range("I88").Value = 2849.8333
$ws.Range("J88").Value = 3125
$ws.Range("K88").Value = 2849.8333
$ws.Range("L88").Value = 3125
$ws.Range("M88").Value = -2443.8333
$ws.Range("N88").Value = -3937

$ws.Range("H91").Value = 2959.9
$ws.Range("I91").Value = 2849.8333
$ws.Range("J91").Value = 3125
$ws.Range("K91").Value = 2849.8333
$ws.Range("L91").Value = 3125
$ws.Range("M91").Value = -1445.8333
$ws.Range("N91").Value = -5933

$ws.Range("H123").Value = 38097.832
$ws.Range("J123").Value = 38097.832
$ws.Range("L123").Value = 38097.832
$ws.Range("N123").Value = -47897.832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 17859422
$ws.Range("I86").Value = 2150.6667
$ws.Range("J86").Value = 50002510
$ws.Range("K86").Value = 2150.6667
$ws.Range("L86").Value = 50002510
$ws.Range("M86").Value = -1027.6667
$ws.Range("N86").Value = -50004756

$ws.Range("H89").Value = 17859422
$ws.Range("I89").Value = 2150.6667
$ws.Range("J89").Value = 50002510
$ws.Range("K89").Value = 10753.3335
$ws.Range("L89").Value = 250012550
$ws.Range("M89").Value = -5137.333500000001
$ws.Range("N89").Value = -250023782

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 14708476
$ws.Range("I132").Value = 17858970
$ws.Range("J132").Value = 6171
$ws.Range("K132").Value = 53576910
$ws.Range("L132").Value = 18513
$ws.Range("M132").Value = -53574380
$ws.Range("N132").Value = -23573

$ws.Range("H140").Value = 40158.383
$ws.Range("J140").Value = 40158.383
$ws.Range("L140").Value = 40158.383
$ws.Range("N140").Value = -50518.383

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4999.2144
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4999.2144
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 14997.6432
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -15627.6432

$ws.Range("H73").Value = 4999.2144
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4999.2144
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 14997.6432
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -17181.6432

$ws.Range("H87").Value = 11883.333
$ws.Range("I87").Value = 5540
$ws.Range("J87").Value = 19812.5
$ws.Range("K87").Value = 16620
$ws.Range("L87").Value = 59437.5
$ws.Range("M87").Value = -15372
$ws.Range("N87").Value = -61933.5

$ws.Range("H90").Value = 11883.333
$ws.Range("I90").Value = 5540
$ws.Range("J90").Value = 19812.5
$ws.Range("K90").Value = 49860
$ws.Range("L90").Value = 178312.5
$ws.Range("M90").Value = -43620
$ws.Range("N90").Value = -190792.5

$ws.Range("H117").Value = 595.6
$ws.Range("J117").Value = 595.6
$ws.Range("L117").Value = 1786.8
$ws.Range("N117").Value = -8670.799999999999

$ws.Range("H129").Value = 2652.9644
$ws.Range("I129").Value = 1166.6666
$ws.Range("J129").Value = 3357
$ws.Range("K129").Value = 3499.9998
$ws.Range("L129").Value = 10071
$ws.Range("M129").Value = 1500.0002
$ws.Range("N129").Value = -20071

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 18088.5
$ws.Range("J15").Value = 18088.5
$ws.Range("L15").Value = 18088.5
$ws.Range("N15").Value = -18664.5

$ws.Range("H63").Value = 39900
$ws.Range("J63").Value = 39900
$ws.Range("L63").Value = 39900
$ws.Range("N63").Value = -41272

$ws.Range("H66").Value = 39900
$ws.Range("J66").Value = 39900
$ws.Range("L66").Value = 119700
$ws.Range("N66").Value = -126564

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H80").Value = 17547568
$ws.Range("I80").Value = 23811950
$ws.Range("J80").Value = 7300
$ws.Range("K80").Value = 23811950
$ws.Range("L80").Value = 7300
$ws.Range("M80").Value = -23810952
$ws.Range("N80").Value = -9296

$ws.Range("H81").Value = 18088.5
$ws.Range("J81").Value = 18088.5
$ws.Range("L81").Value = 18088.5
$ws.Range("N81").Value = -20084.5

$ws.Range("H83").Value = 17547568
$ws.Range("I83").Value = 23811950
$ws.Range("J83").Value = 7300
$ws.Range("K83").Value = 119059750
$ws.Range("L83").Value = 36500
$ws.Range("M83").Value = -119054758
$ws.Range("N83").Value = -46484

$ws.Range("H84").Value = 18088.5
$ws.Range("J84").Value = 18088.5
$ws.Range("L84").Value = 54265.5
$ws.Range("N84").Value = -64249.5

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H102").Value = 4022.0952
$ws.Range("I102").Value = 5143.2144
$ws.Range("J102").Value = 1779.8572
$ws.Range("K102").Value = 5143.2144
$ws.Range("L102").Value = 1779.8572
$ws.Range("M102").Value = -3521.2144
$ws.Range("N102").Value = -5023.8572

$ws.Range("H106").Value = 32940
$ws.Range("J106").Value = 32940
$ws.Range("L106").Value = 32940
$ws.Range("N106").Value = -35464

$ws.Range("H113").Value = 68026.664
$ws.Range("I113").Value = 126126.25
$ws.Range("J113").Value = 1627.1428
$ws.Range("K113").Value = 126126.25
$ws.Range("L113").Value = 1627.1428
$ws.Range("M113").Value = -123956.25
$ws.Range("N113").Value = -5967.1428

$ws.Range("H126").Value = 3482.1042
$ws.Range("I126").Value = 2630.0356
$ws.Range("J126").Value = 4675
$ws.Range("K126").Value = 7890.1068
$ws.Range("L126").Value = 14025
$ws.Range("M126").Value = -5420.1068
$ws.Range("N126").Value = -18965

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 910.1
$ws.Range("I22").Value = 640.2
$ws.Range("J22").Value = 1180
$ws.Range("K22").Value = 640.2
$ws.Range("L22").Value = 1180
$ws.Range("M22").Value = -345.2
$ws.Range("N22").Value = -1770

$ws.Range("H27").Value = 910.1
$ws.Range("I27").Value = 640.2
$ws.Range("J27").Value = 1180
$ws.Range("K27").Value = 640.2
$ws.Range("L27").Value = 1180
$ws.Range("M27").Value = -533.2
$ws.Range("N27").Value = -1394

$ws.Range("H122").Value = 6686.6
$ws.Range("I122").Value = 7469.3335
$ws.Range("J122").Value = 5512.5
$ws.Range("K122").Value = 22408.0005
$ws.Range("L122").Value = 16537.5
$ws.Range("M122").Value = -19958.0005
$ws.Range("N122").Value = -21437.5

$ws.Range("H139").Value = 59818
$ws.Range("J139").Value = 59818
$ws.Range("L139").Value = 59818
$ws.Range("N139").Value = -70098

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11130.272
$ws.Range("I62").Value = 4664.2856
$ws.Range("J62").Value = 22445.75
$ws.Range("K62").Value = 4664.2856
$ws.Range("L62").Value = 22445.75
$ws.Range("M62").Value = -4040.2856
$ws.Range("N62").Value = -23693.75

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H65").Value = 11130.272
$ws.Range("I65").Value = 4664.2856
$ws.Range("J65").Value = 22445.75
$ws.Range("K65").Value = 23321.428
$ws.Range("L65").Value = 112228.75
$ws.Range("M65").Value = -20201.428
$ws.Range("N65").Value = -118468.75

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H70").Value = 20102.5
$ws.Range("J70").Value = 20102.5
$ws.Range("L70").Value = 20102.5
$ws.Range("N70").Value = -20732.5

$ws.Range("H73").Value = 20102.5
$ws.Range("J73").Value = 20102.5
$ws.Range("L73").Value = 20102.5
$ws.Range("N73").Value = -22286.5
